# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets
# to match the newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) updates for column F
$updates = @{
    "展览" = @(
        @{Row = 9;  Value = 1144},
        @{Row = 10; Value = 3320},
        @{Row = 11; Value = 2456},
        @{Row = 13; Value = 2251},
        @{Row = 18; Value = 610},
        @{Row = 20; Value = 272},
        @{Row = 26; Value = 16},
        @{Row = 28; Value = 432},
        @{Row = 29; Value = 627},
        @{Row = 35; Value = 634},
        @{Row = 36; Value = 629},
        @{Row = 37; Value = 1808},
        @{Row = 39; Value = 476},
        @{Row = 41; Value = 504},
        @{Row = 42; Value = 1122},
        @{Row = 44; Value = 384}
    )
    "全部类型" = @(
        @{Row = 8;  Value = 1144},
        @{Row = 9;  Value = 3320},
        @{Row = 10; Value = 2456},
        @{Row = 11; Value = 2251},
        @{Row = 16; Value = 610},
        @{Row = 18; Value = 272},
        @{Row = 24; Value = 16},
        @{Row = 26; Value = 432},
        @{Row = 27; Value = 627},
        @{Row = 36; Value = 634},
        @{Row = 38; Value = 629},
        @{Row = 39; Value = 1808},
        @{Row = 44; Value = 476},
        @{Row = 46; Value = 504},
        @{Row = 47; Value = 1122},
        @{Row = 48; Value = 384}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
